$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")

# "Research and Upgrade => MAC + Win 8 + VS2013" (row 8) moved from "ALMOST DONE"
# to "IN PROGRESS" status.
$ws.Range("B8").Value = "IN PROGRESS"

# Insert a brand-new row right below it (native full-row insert keeps each
# shifted cell's exact style index intact).
$ws.Rows.Item(9).Insert(-4121)  # -4121 = xlShiftDown

# The only cell below the insertion point outside columns A:B was the blank
# divider cell C10; a full-row insert also pushed it down to C11, but the
# real edit only shifted columns A:B. Put the divider back at C10 and drop
# the now-stray C11 cell entirely.
$ws.Range("C10").HorizontalAlignment = -4152  # xlRight
$ws.Range("C11").Clear()

# Fill in the newly inserted row with the new backlog item.
$ws.Range("A9").Value = "Redbox - login to ADP, etc."
$ws.Range("B9").Value = "TODO"

# Leave the selection the way the author ended up leaving it.
$ws.Activate()
[void]$ws.Range("A7:B16").Select()
